# Add the "Adress" (column C) data for the OffProgram/HallList data model
# and append a new "Инфоцентр" row, per commit "add data model OffProgram
# and upload at data_loader".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column C ("Adress") values for the existing 8 halls --------------
# (written in this order so the shared-string table is populated with the
#  same new unique-string order as the target workbook)
$ws.Range("C8").Value = "ул. К. Либкнехта, 48"
$ws.Range("C3").Value = "ул. Я. Свердлова, 30"
$ws.Range("C6").Value = "ул. Пролетарская, 18"
$ws.Range("C2").Value = "ул. К. Либкнехта, 38а"
$ws.Range("C5").Value = "ул. К. Либкнехта, 38а"
$ws.Range("C4").Value = "ул. Первомайская, 22"
$ws.Range("C7").Value = "ул. Первомайская, 9"
$ws.Range("C9").Value = "ул. К. Либкнехта, 48"

# --- New row 10: "Инфоцентр" ------------------------------------------------
$ws.Range("A10").Value = "Инфоцентр"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = "ул. Первомайская, 9"
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0

# Give A10 the same header-ish look as A1:A9 (bold font, centered/top
# aligned) but with only left/right borders instead of a full box.
$ws.Range("A1").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Borders.Item(8).LineStyle = -4142
$ws.Range("A10").Borders.Item(9).LineStyle = -4142

# --- Column C width ----------------------------------------------------
$ws.Range("C1").ColumnWidth = 21.83

# --- Selection moves to A11 after the new last row ----------------------
$ws.Range("A11").Select() | Out-Null

Write-Output "edit applied"
